$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original style of the Price column while writing text-like numeric
# strings (e.g. "29.180.03", "1.004", "0.000009642") so Excel does not
# auto-convert them to numbers and strip formatting/zeros.
$priceRange = $ws.Range("D2:D51")
$savedStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value2 = "29.180.03"
$ws.Range("E2").Value2 = "  -0.66%  "

$ws.Range("D3").Value2 = "1.826.45"
$ws.Range("E3").Value2 = "  -0.79%  "

$ws.Range("D4").Value2 = "1.004"
$ws.Range("E4").Value2 = "  +0.43%  "

$ws.Range("D5").Value2 = "233.54"
$ws.Range("E5").Value2 = "  -2.33%  "

$ws.Range("D6").Value2 = "0.5955"
$ws.Range("E6").Value2 = "  -5.08%  "

$ws.Range("D7").Value2 = "1.004"
$ws.Range("E7").Value2 = "  +0.33%  "

$ws.Range("D8").Value2 = "0.06918"
$ws.Range("E8").Value2 = "  -6.45%  "

$ws.Range("D9").Value2 = "0.2731"
$ws.Range("E9").Value2 = "  -5.54%  "

$ws.Range("D10").Value2 = "23.14"
$ws.Range("E10").Value2 = "  -6.85%  "

$ws.Range("D11").Value2 = "0.07594"
$ws.Range("E11").Value2 = "  -1.49%  "

$ws.Range("D12").Value2 = "1.830.62"
$ws.Range("E12").Value2 = "  -0.36%  "

$ws.Range("D13").Value2 = "4.731"
$ws.Range("E13").Value2 = "  -4.72%  "

$ws.Range("D14").Value2 = "0.6211"
$ws.Range("E14").Value2 = "  -7.70%  "

$ws.Range("D15").Value2 = "0.000009642"
$ws.Range("E15").Value2 = "  -5.55%  "

$ws.Range("D16").Value2 = "78.09"
$ws.Range("E16").Value2 = "  -4.53%  "

$ws.Range("D17").Value2 = "28.886.18"
$ws.Range("E17").Value2 = "  -1.58%  "

$ws.Range("D18").Value2 = "5.690"
$ws.Range("E18").Value2 = "  -9.32%  "

$ws.Range("D19").Value2 = "220.67"
$ws.Range("E19").Value2 = "  -5.72%  "

$ws.Range("D20").Value2 = "1.004"
$ws.Range("E20").Value2 = "  +0.35%  "

$ws.Range("D21").Value2 = "11.46"
$ws.Range("E21").Value2 = "  -7.02%  "

$ws.Range("D22").Value2 = "6.831"
$ws.Range("E22").Value2 = "  -6.35%  "

$ws.Range("D23").Value2 = "1.006"
$ws.Range("E23").Value2 = "  +0.51%  "

$ws.Range("D24").Value2 = "155.59"
$ws.Range("E24").Value2 = "  -1.31%  "

$ws.Range("D25").Value2 = "7.900"
$ws.Range("E25").Value2 = "  -6.97%  "

$ws.Range("D26").Value2 = "0.1282"
$ws.Range("E26").Value2 = "  -4.42%  "

$ws.Range("D27").Value2 = "16.43"
$ws.Range("E27").Value2 = "  -5.03%  "

$ws.Range("D28").Value2 = "0.06564"
$ws.Range("E28").Value2 = "  -9.16%  "

$ws.Range("D29").Value2 = "1.447"
$ws.Range("E29").Value2 = "  -3.14%  "

$ws.Range("D30").Value2 = "1.436"
$ws.Range("E30").Value2 = "  -2.69%  "

$ws.Range("D31").Value2 = "3.817"
$ws.Range("E31").Value2 = "  -5.27%  "

$ws.Range("D32").Value2 = "3.736"
$ws.Range("E32").Value2 = "  -7.47%  "

$ws.Range("D33").Value2 = "1.085"
$ws.Range("E33").Value2 = "  -5.50%  "

$ws.Range("D34").Value2 = "1.705"
$ws.Range("E34").Value2 = "  -6.08%  "

$ws.Range("D35").Value2 = "0.6381"
$ws.Range("E35").Value2 = "  -8.54%  "

$ws.Range("D36").Value2 = "2.542"
$ws.Range("E36").Value2 = "  -1.35%  "

$ws.Range("D37").Value2 = "2.734"
$ws.Range("E37").Value2 = "  -2.41%  "

$ws.Range("B38").Value2 = "VeChain"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value2 = "0.01726"
$ws.Range("E38").Value2 = "  -5.61%  "

$ws.Range("B39").Value2 = "FraxShare"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value2 = "6.493"
$ws.Range("E39").Value2 = "  -4.05%  "

$ws.Range("B40").Value2 = "Maker"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value2 = "1.173.25"
$ws.Range("E40").Value2 = "  -4.80%  "

$ws.Range("D41").Value2 = "0.8942"
$ws.Range("E41").Value2 = "  -5.52%  "

$ws.Range("D42").Value2 = "1.004"
$ws.Range("E42").Value2 = "  +0.29%  "

$ws.Range("D43").Value2 = "1.979.74"
$ws.Range("E43").Value2 = "  -0.74%  "

$ws.Range("D44").Value2 = "100.26"
$ws.Range("E44").Value2 = "  -0.84%  "

$ws.Range("D45").Value2 = "61.64"
$ws.Range("E45").Value2 = "  -5.58%  "

$ws.Range("D46").Value2 = "0.00000000115"
$ws.Range("E46").Value2 = "  -3.75%  "

$ws.Range("B47").Value2 = "Cronos"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value2 = "0.05540"
$ws.Range("E47").Value2 = "  -2.21%  "

$ws.Range("B48").Value2 = "EnergySwap"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value2 = "8.419"
$ws.Range("E48").Value2 = "  -5.45%  "

$ws.Range("D49").Value2 = "0.4555"
$ws.Range("E49").Value2 = "  -0.47%  "

$ws.Range("D50").Value2 = "1.562"
$ws.Range("E50").Value2 = "  -7.90%  "

$ws.Range("B51").Value2 = "TheSandbox"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value2 = "0.3611"
$ws.Range("E51").Value2 = "  -7.09%  "

# Restore the original (default) style/number format on the Price column.
$priceRange.Style = $savedStyle
